$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 323.5
$ws.Range("I28").Value = 338.33334
$ws.Range("J28").Value = 190
$ws.Range("K28").Value = 338.33334
$ws.Range("L28").Value = 190
$ws.Range("M28").Value = 146.66666
$ws.Range("N28").Value = -1160
$ws.Range("H107").Value = 651.2692
$ws.Range("I107").Value = 556.55
$ws.Range("J107").Value = 967
$ws.Range("K107").Value = 556.55
$ws.Range("L107").Value = 967
$ws.Range("M107").Value = 1363.45
$ws.Range("N107").Value = -4807
$ws.Range("H111").Value = 1514.5
$ws.Range("I111").Value = 1529
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 4587
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -1520
$ws.Range("N111").Value = -10634
$ws.Range("H129").Value = 808.14
$ws.Range("I129").Value = 473.8
$ws.Range("J129").Value = 891.725
$ws.Range("K129").Value = 1421.4
$ws.Range("L129").Value = 2675.175
$ws.Range("M129").Value = 3578.6
$ws.Range("N129").Value = -12675.175

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3485.9858
$ws.Range("I32").Value = 1899.25
$ws.Range("J32").Value = 12140.909
$ws.Range("K32").Value = 1899.25
$ws.Range("L32").Value = 12140.909
$ws.Range("M32").Value = -1612.25
$ws.Range("N32").Value = -12714.909
$ws.Range("H61").Value = 2423.027
$ws.Range("I61").Value = 1864
$ws.Range("J61").Value = 3242.9333
$ws.Range("K61").Value = 1864
$ws.Range("L61").Value = 3242.9333
$ws.Range("M61").Value = -1652
$ws.Range("N61").Value = -3666.9333
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H110").Value = 1651.6471
$ws.Range("I110").Value = 1651.6471
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1651.6471
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 393.3529000000001
$ws.Range("N110").Value = $null
$ws.Range("H112").Value = 42686.5
$ws.Range("J112").Value = 42686.5
$ws.Range("L112").Value = 42686.5
$ws.Range("N112").Value = -45640.5
$ws.Range("H122").Value = 1930.0588
$ws.Range("I122").Value = 1807.9286
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5423.7858
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2973.7858
$ws.Range("N122").Value = -12400
$ws.Range("H136").Value = 2423.027
$ws.Range("I136").Value = 1864
$ws.Range("J136").Value = 3242.9333
$ws.Range("K136").Value = 5592
$ws.Range("L136").Value = 9728.7999
$ws.Range("M136").Value = -3042
$ws.Range("N136").Value = -14828.7999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H93").Value = 18655.555
$ws.Range("J93").Value = 18655.555
$ws.Range("L93").Value = 18655.555
$ws.Range("N93").Value = -22399.555
$ws.Range("H97").Value = 15747.777
$ws.Range("I97").Value = 865
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 865
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = 126
$ws.Range("N97").Value = -21982
$ws.Range("H107").Value = 1529.6316
$ws.Range("I107").Value = 1060.7858
$ws.Range("J107").Value = 2842.4
$ws.Range("K107").Value = 1060.7858
$ws.Range("L107").Value = 2842.4
$ws.Range("M107").Value = 859.2141999999999
$ws.Range("N107").Value = -6682.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1130.2727
$ws.Range("I16").Value = 1159.2222
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1159.2222
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -872.2221999999999
$ws.Range("N16").Value = -1574
$ws.Range("H58").Value = 1889.4872
$ws.Range("I58").Value = 1246.963
$ws.Range("J58").Value = 3335.1667
$ws.Range("K58").Value = 1246.963
$ws.Range("L58").Value = 3335.1667
$ws.Range("M58").Value = -1043.963
$ws.Range("N58").Value = -3741.1667
$ws.Range("H76").Value = 2666.5
$ws.Range("I76").Value = 2666.5
$ws.Range("K76").Value = 2666.5
$ws.Range("M76").Value = -2351.5
$ws.Range("H79").Value = 2666.5
$ws.Range("I79").Value = 2666.5
$ws.Range("K79").Value = 2666.5
$ws.Range("M79").Value = -1574.5
$ws.Range("H107").Value = 775
$ws.Range("I107").Value = 614.4286
$ws.Range("K107").Value = 614.4286
$ws.Range("M107").Value = 1305.5714
$ws.Range("H113").Value = 1130.2727
$ws.Range("I113").Value = 1159.2222
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1159.2222
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1010.7778
$ws.Range("N113").Value = -5340
$ws.Range("H134").Value = 388548.28
$ws.Range("I134").Value = 403641
$ws.Range("J134").Value = 282899.2
$ws.Range("K134").Value = 1210923
$ws.Range("L134").Value = 848697.6000000001
$ws.Range("M134").Value = -1208388
$ws.Range("N134").Value = -853767.6000000001
$ws.Range("H136").Value = 1889.4872
$ws.Range("I136").Value = 1246.963
$ws.Range("J136").Value = 3335.1667
$ws.Range("K136").Value = 3740.889
$ws.Range("L136").Value = 10005.5001
$ws.Range("M136").Value = -1190.889
$ws.Range("N136").Value = -15105.5001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 831.7692
$ws.Range("I8").Value = 831.7692
$ws.Range("K8").Value = 2495.3076
$ws.Range("M8").Value = -2356.3076
$ws.Range("H98").Value = 1130.5
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 1256.6
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 3769.8
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -6765.799999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5640.409
$ws.Range("I80").Value = 4005.9333
$ws.Range("J80").Value = 9142.857
$ws.Range("K80").Value = 4005.9333
$ws.Range("L80").Value = 9142.857
$ws.Range("M80").Value = -3007.9333
$ws.Range("N80").Value = -11138.857
$ws.Range("H83").Value = 5640.409
$ws.Range("I83").Value = 4005.9333
$ws.Range("J83").Value = 9142.857
$ws.Range("K83").Value = 20029.6665
$ws.Range("L83").Value = 45714.285
$ws.Range("M83").Value = -15037.6665
$ws.Range("N83").Value = -55698.285
$ws.Range("H107").Value = 6233.3335
$ws.Range("I107").Value = 150
$ws.Range("J107").Value = 7450
$ws.Range("K107").Value = 150
$ws.Range("L107").Value = 7450
$ws.Range("M107").Value = 1770
$ws.Range("N107").Value = -11290
$ws.Range("H122").Value = 1685.7142
$ws.Range("I122").Value = 1760
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5280
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2830
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 45457580
$ws.Range("I132").Value = 76925270
$ws.Range("J132").Value = 4243.222
$ws.Range("K132").Value = 230775810
$ws.Range("L132").Value = 12729.666
$ws.Range("M132").Value = -230773280
$ws.Range("N132").Value = -17789.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2713.0625
$ws.Range("I7").Value = 2075.3333
$ws.Range("J7").Value = 4626.25
$ws.Range("K7").Value = 2075.3333
$ws.Range("L7").Value = 4626.25
$ws.Range("M7").Value = -1963.3333
$ws.Range("N7").Value = -4850.25
$ws.Range("H61").Value = 3016.125
$ws.Range("I61").Value = 3662.25
$ws.Range("J61").Value = 2370
$ws.Range("K61").Value = 3662.25
$ws.Range("L61").Value = 2370
$ws.Range("M61").Value = -3460.25
$ws.Range("N61").Value = -2774
$ws.Range("H113").Value = 3016.125
$ws.Range("I113").Value = 3662.25
$ws.Range("J113").Value = 2370
$ws.Range("K113").Value = 3662.25
$ws.Range("L113").Value = 2370
$ws.Range("M113").Value = -1492.25
$ws.Range("N113").Value = -6710
$ws.Range("H122").Value = 85398.164
$ws.Range("I122").Value = 144871.14
$ws.Range("J122").Value = 2136
$ws.Range("K122").Value = 434613.42
$ws.Range("L122").Value = 6408
$ws.Range("M122").Value = -432163.42
$ws.Range("N122").Value = -11308
$ws.Range("H126").Value = 2713.0625
$ws.Range("I126").Value = 2075.3333
$ws.Range("J126").Value = 4626.25
$ws.Range("K126").Value = 6225.999899999999
$ws.Range("L126").Value = 13878.75
$ws.Range("M126").Value = -3755.999899999999
$ws.Range("N126").Value = -18818.75
$ws.Range("H136").Value = 2393.8823
$ws.Range("I136").Value = 1508.5
$ws.Range("J136").Value = 4518.8
$ws.Range("K136").Value = 4525.5
$ws.Range("L136").Value = 13556.4
$ws.Range("M136").Value = -1975.5
$ws.Range("N136").Value = -18656.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 710.375
$ws.Range("I113").Value = 710.375
$ws.Range("K113").Value = 2131.125
$ws.Range("M113").Value = 38.875
$ws.Range("H122").Value = 2599278.8
$ws.Range("I122").Value = 2599278.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7797836.399999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7795386.399999999
$ws.Range("N122").Value = $null
$ws.Range("H136").Value = 1112342.9
$ws.Range("I136").Value = 1556534.2
$ws.Range("K136").Value = 4669602.6
$ws.Range("M136").Value = -4667052.6
